$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C3 and C4 from "Invalid" to "Valid"
$ws.Range("C3").Value = "Valid"
$ws.Range("C4").Value = "Valid"

# Update the selected cell / range shown in the sheet view
$ws.Range("F4").Select()
